$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '37.810.33'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +0.10%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.082.92'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +0.30%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '233.71'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.37%  '
$ws.Range("E6").Value = '  -0.09%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '58.97'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +3.07%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("E9").Value = '  +1.97%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0789'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +0.09%  '
$ws.Range("E11").Value = '  +2.80%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.389.06'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +0.58%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.78'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +2.32%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '21.27'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +1.45%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.775'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +1.74%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.33'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +1.25%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.070.15'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -0.27%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '37.725.16'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +0.13%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.15'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +0.41%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '71.80'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +1.70%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0847'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +3.18%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '228.39'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.18%  '
$ws.Range("E23").Value = '  -0.09%  '
$ws.Range("E24").Value = '  -0.48%  '
$ws.Range("E25").Value = '  +0.90%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.63'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +8.10%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '171.16'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +0.57%  '
$ws.Range("E28").Value = '  -1.48%  '
$ws.Range("E29").Value = '  -1.28%  '
$ws.Range("E30").Value = '  +0.93%  '
$ws.Range("E31").Value = '  +2.23%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.74'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +2.19%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0635'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +1.75%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.67'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +1.23%  '
$ws.Range("E35").Value = '  -0.28%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.44'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +1.13%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.83'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +0.01%  '
$ws.Range("E38").Value = '  +0.19%  '
$ws.Range("E39").Value = '  +0.45%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0982'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -1.87%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '17.35'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +11.24%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '99.02'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +1.00%  '
$ws.Range("E43").Value = '  +2.69%  '
$ws.Range("E44").Value = '  -1.35%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.450.86'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +0.17%  '
$ws.Range("E46").Value = '  -0.56%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.16'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +0.84%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.06'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +1.38%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.36'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -0.76%  '
$ws.Range("E50").Value = '  -0.70%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.275.30'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +0.37%  '
